$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add the new data row: A4 = 777777 (extends used range to A1:A4)
$ws.Range("A4").Value = 777777

# Match the saved cursor position recorded in the file (C9)
$ws.Range("C9").Select()
